$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# 1. Insert two rows at the top. Everything that used to be on rows 1-14
#    shifts down to rows 3-16. Row 1 stays empty, row 2 becomes the new
#    title row.
# -----------------------------------------------------------------------
$ws.Rows("1:2").Insert()

# -----------------------------------------------------------------------
# 2. Title row (row 2)
# -----------------------------------------------------------------------
$ws.Range("D2").Value = "Permissions for People Using the Missionout Software"
$ws.Range("D2").Font.Bold = $true
$ws.Range("D2").Font.Italic = $true
$ws.Range("D2").Font.Size = 24
$ws.Rows("2:2").RowHeight = 31

# -----------------------------------------------------------------------
# 3. Bump every font in the main table (rows 4-10, the old rows 2-8) up
#    to match the new, larger type sizes used throughout the chart.
# -----------------------------------------------------------------------
$ws.Range("B4:J5").Font.Size = 18
$ws.Range("B6:J10").Font.Size = 18

# -----------------------------------------------------------------------
# 4. New row 3 (formerly the blank divider row 1) and the new column A
#    border/format column that now frames the table on rows 3-10.
# -----------------------------------------------------------------------
$ws.Range("A3,D3:J3").Font.Size = 16
$ws.Range("C3").Font.Size = 16
$ws.Range("C3").Font.Italic = $true

$ws.Range("A5").Font.Size = 16
$ws.Range("A5").Font.Italic = $true
$ws.Range("A6:A10").Font.Size = 16

# -----------------------------------------------------------------------
# 5. Clear the old inner grid-lines on the data rows (Editors/Users/
#    Signed In no longer carry the left/right "medium" rules down each
#    column in the body of the table).
# -----------------------------------------------------------------------
$ws.Range("D6:I8").Borders.Item(7).LineStyle = -4142
$ws.Range("D6:I8").Borders.Item(10).LineStyle = -4142
$ws.Range("D4:I4").Borders.Item(7).LineStyle = -4142
$ws.Range("D4:I4").Borders.Item(10).LineStyle = -4142
$ws.Range("D5:I5").Borders.Item(7).LineStyle = -4142
$ws.Range("D5:I5").Borders.Item(10).LineStyle = -4142

# Header row (row4) columns D:I keep only the top "medium" rule and
# become center-aligned.
$ws.Range("D4:I4").HorizontalAlignment = -4108

# -----------------------------------------------------------------------
# 6. Give the Anyomous row (row 9) its own thin box instead of the
#    medium one that used to close off the table.
# -----------------------------------------------------------------------
$ws.Range("D9:I9").Borders.Item(9).LineStyle = 1
$ws.Range("D9:I9").Borders.Item(9).Weight = 2
$ws.Range("D9").Borders.Item(7).LineStyle = 1
$ws.Range("D9").Borders.Item(7).Weight = 2
$ws.Range("I9").Borders.Item(10).LineStyle = 1
$ws.Range("I9").Borders.Item(10).Weight = 2

# -----------------------------------------------------------------------
# 7. Fill in the permissions grid (the actual "chart").
# -----------------------------------------------------------------------
$ws.Range("D6:I6").Value = "Write/Create"
$ws.Range("D7:I7").Value = "Write/Create"
$ws.Range("F7").Value = "Read Access"
$ws.Range("D8").Value = "Write/Create"
$ws.Range("E8:I8").Value = "Read Access"
$ws.Range("D9:I9").Value = "No Access"

# -----------------------------------------------------------------------
# 8. Echo columns: J6:J9 repeat the row label, D10:I10 repeat the
#    column header - both done with formulas so they recalc live.
# -----------------------------------------------------------------------
$ws.Range("J6:J9").Formula = "=B6"
$ws.Range("D10:I10").Formula = "=D4"

# -----------------------------------------------------------------------
# 9. Footer ("Updated" + date) and the legend restated below it.
# -----------------------------------------------------------------------
$ws.Range("B12").Value = "Updated"
$ws.Range("B12").Font.Size = 16
$ws.Range("B12").HorizontalAlignment = -4152

$ws.Range("B13").Value = 43922
$ws.Range("B13").Font.Size = 16
$ws.Range("B13").HorizontalAlignment = -4152
$ws.Range("B13").NumberFormat = "d-mmm-yy"

$ws.Range("E13").Font.Size = 18
$ws.Range("E13").Font.Bold = $true
$ws.Range("E13").Font.Italic = $true

$ws.Range("F13:F16").Font.Size = 18

# -----------------------------------------------------------------------
# 10. Column widths - let Excel re-fit them against the new (larger)
#     text now that every cell is filled in.
# -----------------------------------------------------------------------
$ws.Columns("C:J").EntireColumn.AutoFit()

$ws.Range("C20").Select()
